$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07285266666666666
$ws.Range("H2").Value = 0.218558
$ws.Range("I2").Value = 0.05584899373277382
$ws.Range("J2").Value = 0.05584899373277381
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 0.05419946989333332
$ws.Range("R2").Value = 0.48779522904
$ws.Range("S2").Value = 0.000349559031428712
$ws.Range("T2").Value = 0.000349559031428712
# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07285266666666666
$ws.Range("H3").Value = 0.218558
$ws.Range("I3").Value = 0.05584899373277382
$ws.Range("J3").Value = 0.05584899373277381
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("Q3").Value = 6.421946927627554
$ws.Range("R3").Value = 57.79752234864799
$ws.Range("S3").Value = 0.04141829343213194
$ws.Range("T3").Value = 0.04141829343213194
# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07285266666666666
$ws.Range("H4").Value = 0.218558
$ws.Range("I4").Value = 0.05584899373277382
$ws.Range("J4").Value = 0.05584899373277381
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 2.168721844364889
$ws.Range("R4").Value = 19.518496599284
$ws.Range("S4").Value = 0.01398715354313324
$ws.Range("T4").Value = 0.01398715354313323
# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07285266666666666
$ws.Range("H5").Value = 0.218558
$ws.Range("I5").Value = 0.05584899373277382
$ws.Range("J5").Value = 0.05584899373277381
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 0.01457288890288889
$ws.Range("R5").Value = 0.131156000126
$ws.Range("S5").Value = 0.00009398772607993074
$ws.Range("T5").Value = 0.00009398772607993074
# Row 6
$ws.Range("I6").Value = 0.438974399073536
$ws.Range("J6").Value = 0.438974399073536
$ws.Range("M6").Value = 0.74396
$ws.Range("N6").Value = 2.23188
$ws.Range("O6").Value = 0.006259003216804254
$ws.Range("P6").Value = 0.006259003216804255
$ws.Range("Q6").Value = 0.4260091030533333
$ws.Range("R6").Value = 3.83408192748
$ws.Range("S6").Value = 0.002747542175895977
$ws.Range("T6").Value = 0.002747542175895977
# Row 7
$ws.Range("I7").Value = 0.438974399073536
$ws.Range("J7").Value = 0.438974399073536
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("O7").Value = 0.7416121699579786
$ws.Range("P7").Value = 0.7416121699579786
$ws.Range("Q7").Value = 50.47665329345288
$ws.Range("R7").Value = 454.2898796410759
$ws.Range("S7").Value = 0.3255487566529247
$ws.Range("T7").Value = 0.3255487566529247
# Row 8
$ws.Range("I8").Value = 0.438974399073536
$ws.Range("J8").Value = 0.438974399073536
$ws.Range("M8").Value = 29.76859933333333
$ws.Range("N8").Value = 89.305798
$ws.Range("O8").Value = 0.2504459365921425
$ws.Range("P8").Value = 0.2504459365921425
$ws.Range("Q8").Value = 17.04620450178422
$ws.Range("R8").Value = 153.415840516058
$ws.Range("S8").Value = 0.1099393545159447
$ws.Range("T8").Value = 0.1099393545159447
# Row 9
$ws.Range("I9").Value = 0.438974399073536
$ws.Range("J9").Value = 0.438974399073536
$ws.Range("M9").Value = 0.2000323333333334
$ws.Range("N9").Value = 0.6000970000000001
$ws.Range("O9").Value = 0.00168289023307462
$ws.Range("P9").Value = 0.00168289023307462
$ws.Range("Q9").Value = 0.1145432481652222
$ws.Range("R9").Value = 1.030889233487
$ws.Range("S9").Value = 0.0007387457287706544
$ws.Range("T9").Value = 0.0007387457287706545
# Row 10
$ws.Range("G10").Value = 0.594248
$ws.Range("H10").Value = 1.782744
$ws.Range("I10").Value = 0.4555516544035914
$ws.Range("J10").Value = 0.4555516544035914
$ws.Range("M10").Value = 0.74396
$ws.Range("N10").Value = 2.23188
$ws.Range("O10").Value = 0.006259003216804254
$ws.Range("P10").Value = 0.006259003216804255
$ws.Range("Q10").Value = 0.44209674208
$ws.Range("R10").Value = 3.97887067872
$ws.Range("S10").Value = 0.002851299270332579
$ws.Range("T10").Value = 0.002851299270332579
# Row 11
$ws.Range("G11").Value = 0.594248
$ws.Range("H11").Value = 1.782744
$ws.Range("I11").Value = 0.4555516544035914
$ws.Range("J11").Value = 0.4555516544035914
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("O11").Value = 0.7416121699579786
$ws.Range("P11").Value = 0.7416121699579786
$ws.Range("Q11").Value = 52.38283363476266
$ws.Range("R11").Value = 471.445502712864
$ws.Range("S11").Value = 0.3378426509501946
$ws.Range("T11").Value = 0.3378426509501946
# Row 12
$ws.Range("G12").Value = 0.594248
$ws.Range("H12").Value = 1.782744
$ws.Range("I12").Value = 0.4555516544035914
$ws.Range("J12").Value = 0.4555516544035914
$ws.Range("M12").Value = 29.76859933333333
$ws.Range("N12").Value = 89.305798
$ws.Range("O12").Value = 0.2504459365921425
$ws.Range("P12").Value = 0.2504459365921425
$ws.Range("Q12").Value = 17.68993061663467
$ws.Range("R12").Value = 159.209375549712
$ws.Range("S12").Value = 0.1140910607532075
$ws.Range("T12").Value = 0.1140910607532075
# Row 13
$ws.Range("G13").Value = 0.594248
$ws.Range("H13").Value = 1.782744
$ws.Range("I13").Value = 0.4555516544035914
$ws.Range("J13").Value = 0.4555516544035914
$ws.Range("M13").Value = 0.2000323333333334
$ws.Range("N13").Value = 0.6000970000000001
$ws.Range("O13").Value = 0.00168289023307462
$ws.Range("P13").Value = 0.00168289023307462
$ws.Range("Q13").Value = 0.1188688140186667
$ws.Range("R13").Value = 1.069819326168
$ws.Range("S13").Value = 0.0007666434298567887
$ws.Range("T13").Value = 0.0007666434298567888
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.06473366666666668
$ws.Range("H14").Value = 0.194201
$ws.Range("I14").Value = 0.04962495279009878
$ws.Range("J14").Value = 0.04962495279009878
$ws.Range("M14").Value = 0.74396
$ws.Range("N14").Value = 2.23188
$ws.Range("O14").Value = 0.006259003216804254
$ws.Range("P14").Value = 0.006259003216804255
$ws.Range("Q14").Value = 0.04815925865333334
$ws.Range("R14").Value = 0.43343332788
$ws.Range("S14").Value = 0.0003106027391469875
$ws.Range("T14").Value = 0.0003106027391469875
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.06473366666666668
$ws.Range("H15").Value = 0.194201
$ws.Range("I15").Value = 0.04962495279009878
$ws.Range("J15").Value = 0.04962495279009878
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("O15").Value = 0.7416121699579786
$ws.Range("P15").Value = 0.7416121699579786
$ws.Range("Q15").Value = 5.706258820506222
$ws.Range("R15").Value = 51.356329384556
$ws.Range("S15").Value = 0.03680246892272741
$ws.Range("T15").Value = 0.0368024689227274
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.06473366666666668
$ws.Range("H16").Value = 0.194201
$ws.Range("I16").Value = 0.04962495279009878
$ws.Range("J16").Value = 0.04962495279009878
$ws.Range("M16").Value = 29.76859933333333
$ws.Range("N16").Value = 89.305798
$ws.Range("O16").Value = 0.2504459365921425
$ws.Range("P16").Value = 0.2504459365921425
$ws.Range("Q16").Value = 1.927030586377556
$ws.Range("R16").Value = 17.343275277398
$ws.Range("S16").Value = 0.01242836777985715
$ws.Range("T16").Value = 0.01242836777985714
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.06473366666666668
$ws.Range("H17").Value = 0.194201
$ws.Range("I17").Value = 0.04962495279009878
$ws.Range("J17").Value = 0.04962495279009878
$ws.Range("M17").Value = 0.2000323333333334
$ws.Range("N17").Value = 0.6000970000000001
$ws.Range("O17").Value = 0.00168289023307462
$ws.Range("P17").Value = 0.00168289023307462
$ws.Range("Q17").Value = 0.01294882638855556
$ws.Range("R17").Value = 0.116539437497
$ws.Range("S17").Value = 0.00008351334836724637
$ws.Range("T17").Value = 0.00008351334836724637
